$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column D (shifts existing D:K to E:L)
$ws.Columns("D").Insert()

# Copy number formats from column E into new column D so styles match
$ws.Range("E5:E102").Copy()
$ws.Range("D5:D102").PasteSpecial(-4122)
$ws.Application.CutCopyMode = 0

# Populate new column D with FY2018 values
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 351000
$ws.Range("D9").Value = 69000
$ws.Range("D10").Value = 282000
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 6200
$ws.Range("D15").Value = 167600
$ws.Range("D17").Value = 278200
$ws.Range("D18").Value = 72800
$ws.Range("D20").Value = 72200
$ws.Range("D21").Value = 316800
$ws.Range("D22").Value = 48800
$ws.Range("D23").Value = 96200
$ws.Range("D24").Value = 0
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 96200
$ws.Range("D27").Value = 82400
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = -72200
$ws.Range("D33").Value = 82400
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 82400
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 8000
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 42200
$ws.Range("D44").Value = 0
$ws.Range("D45").Value = 36900
$ws.Range("D46").Value = 0
$ws.Range("D47").Value = 0
$ws.Range("D48").Value = 2649700
$ws.Range("D49").Value = 342000
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 14600
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 3102500
$ws.Range("D57").Value = 45500
$ws.Range("D58").Value = "NA"
$ws.Range("D59").Value = 35900
$ws.Range("D60").Value = 0
$ws.Range("D61").Value = 1325900
$ws.Range("D62").Value = 0
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 1488700
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = -585000
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 1538800
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = 82400
$ws.Range("D83").Value = 171800
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D96").Value = -158900
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = 303800
$ws.Range("D101").Value = 0

# Rows where historical values were also corrected (not a pure shift)
$ws.Range("D70").Value = 75000
$ws.Range("E70").Value = 145000
$ws.Range("F70").Value = 145000
$ws.Range("G70").Value = 139000
$ws.Range("H70").Value = 139000
$ws.Range("I70").Value = 139000
$ws.Range("J70").Value = 69000
$ws.Range("K70").Value = 69000

$ws.Range("D89").Value = 197800
$ws.Range("E89").Value = 162100
$ws.Range("F89").Value = 135800
$ws.Range("G89").Value = 121700
$ws.Range("H89").Value = 96700
$ws.Range("I89").Value = 82700
$ws.Range("J89").Value = 48000
$ws.Range("K89").Value = 17000

$ws.Range("D91").Value = -599400
$ws.Range("E91").Value = -543100
$ws.Range("F91").Value = -408000
$ws.Range("G91").Value = -308300
$ws.Range("H91").Value = -345900
$ws.Range("I91").Value = -261200
$ws.Range("J91").Value = -325800
$ws.Range("K91").Value = -115200

$ws.Range("D94").Value = -507200
$ws.Range("E94").Value = -571600
$ws.Range("F94").Value = -346300
$ws.Range("G94").Value = -372000
$ws.Range("H94").Value = -421700
$ws.Range("I94").Value = -325200
$ws.Range("J94").Value = -417200
$ws.Range("K94").Value = -115000

$ws.Range("D102").Value = -5600
$ws.Range("E102").Value = 6300
$ws.Range("F102").Value = 1400
$ws.Range("G102").Value = -11900
$ws.Range("H102").Value = 17200
$ws.Range("I102").Value = -12300
$ws.Range("J102").Value = 2500
$ws.Range("K102").Value = 14900
